# Updated master department mapping
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Network" -> "NETWORK" (shared string rewritten in place; only rows 34-35 use it)
$rows = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $rows; $r++) {
    if ($ws.Cells.Item($r, 7).Value2 -eq "Network") {
        $ws.Cells.Item($r, 7).Value = "NETWORK"
    }
}

# Fill in the missing "Office" value in column H for rows 31 and 32
$ws.Cells.Item(31, 8).Value = "Office"
$ws.Cells.Item(32, 8).Value = "Office"

# Append four new sessional provider rows (36-39) for the new MSW spine center.
# Columns are populated one whole column at a time (matching how the source
# workbook was authored) so new shared-string entries land in the same order.

# Column A - Department
$ws.Cells.Item(36, 1).Value = "Orthopedics"
$ws.Cells.Item(37, 1).Value = "Orthopedics"
$ws.Cells.Item(38, 1).Value = "Orthopedics"
$ws.Cells.Item(39, 1).Value = "Orthopedics"

# Column B - Department name
$ws.Cells.Item(36, 2).Value = "787 11TH AVE MSW SPINE CENTER"
$ws.Cells.Item(37, 2).Value = "787 11TH AVE MSW SPINE CENTER"
$ws.Cells.Item(38, 2).Value = "787 11TH AVE MSW SPINE CENTER"
$ws.Cells.Item(39, 2).Value = "787 11TH AVE MSW SPINE CENTER"

# Column C - Department ID (matches the existing greyed-out "8005002"-style cells)
$ws.Cells.Item(36, 3).Value = 8792001
$ws.Cells.Item(37, 3).Value = 8792001
$ws.Cells.Item(38, 3).Value = 8792001
$ws.Cells.Item(39, 3).Value = 8792001
$ws.Cells.Item(36, 3).Font.Color = 3355443
$ws.Cells.Item(37, 3).Font.Color = 3355443
$ws.Cells.Item(38, 3).Font.Color = 3355443
$ws.Cells.Item(39, 3).Font.Color = 3355443

# Column D - Provider
$ws.Cells.Item(36, 4).Value = "CHO, SAMUEL K-W"
$ws.Cells.Item(37, 4).Value = "MAHAJER, AMIR"
$ws.Cells.Item(38, 4).Value = "KIM, JUN SUP"
$ws.Cells.Item(39, 4).Value = "BAX, JOSEPH A"

# Column E - NPI
$ws.Cells.Item(36, 5).Value = 1538318209
$ws.Cells.Item(37, 5).Value = 1770829541
$ws.Cells.Item(38, 5).Value = 1962829309
$ws.Cells.Item(39, 5).Value = 1750417283

# Column F - Sessional Managed
$ws.Cells.Item(36, 6).Value = "Department"
$ws.Cells.Item(37, 6).Value = "Department"
$ws.Cells.Item(38, 6).Value = "Department"
$ws.Cells.Item(39, 6).Value = "Department"

# Column G - Sessional Site
$ws.Cells.Item(36, 7).Value = "MSW"
$ws.Cells.Item(37, 7).Value = "MSW"
$ws.Cells.Item(38, 7).Value = "MSW"
$ws.Cells.Item(39, 7).Value = "MSW"

# Column H - Sessional Office
$ws.Cells.Item(36, 8).Value = "Office"
$ws.Cells.Item(37, 8).Value = "Office"
$ws.Cells.Item(38, 8).Value = "Office"
$ws.Cells.Item(39, 8).Value = "Office"

# Selection moves to B11 per the authored workbook state
$ws.Range("B11").Select()
